$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 781, pushing the existing rows 781-836
# down to 783-838 (dimension grows from A1:T836 to A1:T838).
$ws.Rows.Item(781).Resize(2).Insert()

# --- New row 781: Early Glo / Primera / Provincia de Limarí ---
$ws.Cells.Item(781,1).Value = 6
$ws.Cells.Item(781,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(781,3).Value = "Metropolitana"
$ws.Cells.Item(781,4).Value = 44516
$ws.Cells.Item(781,5).Value = 13
$ws.Cells.Item(781,6).Value = "Fruta"
$ws.Cells.Item(781,7).Value = 100103
$ws.Cells.Item(781,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(781,9).Value = 100103006
$ws.Cells.Item(781,10).Value = "Nectarín"
$ws.Cells.Item(781,11).Value = "Early Glo"
$ws.Cells.Item(781,12).Value = "Primera"
$ws.Cells.Item(781,13).Value = 10
$ws.Cells.Item(781,14).Value = 450000
$ws.Cells.Item(781,15).Value = 450000
$ws.Cells.Item(781,16).Value = 450000
$ws.Cells.Item(781,17).Value = "`$/bins (420 kilos)"
$ws.Cells.Item(781,18).Value = "Provincia de Limarí"
$ws.Cells.Item(781,19).Value = 1071
$ws.Cells.Item(781,20).Value = 420

# --- New row 782: Early Glo / Segunda / Provincia de Limarí ---
$ws.Cells.Item(782,1).Value = 6
$ws.Cells.Item(782,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(782,3).Value = "Metropolitana"
$ws.Cells.Item(782,4).Value = 44516
$ws.Cells.Item(782,5).Value = 13
$ws.Cells.Item(782,6).Value = "Fruta"
$ws.Cells.Item(782,7).Value = 100103
$ws.Cells.Item(782,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(782,9).Value = 100103006
$ws.Cells.Item(782,10).Value = "Nectarín"
$ws.Cells.Item(782,11).Value = "Early Glo"
$ws.Cells.Item(782,12).Value = "Segunda"
$ws.Cells.Item(782,13).Value = 17
$ws.Cells.Item(782,14).Value = 400000
$ws.Cells.Item(782,15).Value = 400000
$ws.Cells.Item(782,16).Value = 400000
$ws.Cells.Item(782,17).Value = "`$/bins (420 kilos)"
$ws.Cells.Item(782,18).Value = "Provincia de Limarí"
$ws.Cells.Item(782,19).Value = 952
$ws.Cells.Item(782,20).Value = 420
